# Update "Guaranteed Dispatch Perc by Elec Source" workbook to the new
# fuel-source breakdown (v3.3.1 of the model).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GDPbES")

# --- Header row -------------------------------------------------------
# Row 1, column A gets a new title; bold + wrap text, and the row is
# taller to accommodate the wrapped text.
$ws.Range("A1").Value = "Guaranteed Dispatch Fraction (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# --- Renamed fuel source labels ---------------------------------------
$ws.Range("A2").Value = "hard coal"
$ws.Range("A6").Value = "onshore wind"

# --- New row 13: lignite (all zero, like a brand new source) ---------
$ws.Range("A13").Value = "lignite"
$ws.Range("B13").Value = 0
$ws.Range("C13:AK13").Formula = '=$B13'

# --- New row 14: offshore wind (all zero) -----------------------------
$ws.Range("A14").Value = "offshore wind"
$ws.Range("B14").Value = 0
$ws.Range("C14:AK14").Formula = '=$B14'

# --- New row 15: crude oil (mirrors row 11, petroleum) ----------------
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

# --- New row 16: heavy or residual fuel oil (mirrors row 11) ---------
$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

# --- New row 17: municipal solid waste (mirrors row 9, biomass) ------
$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"
